$d = $word.ActiveDocument

# 1. "Online Server" -> "Database"
$r1 = $d.Content
$r1.Find.Execute("Online Server", $false, $false, $false, $false, $false, $true, 1, $false, "Database", 2)

# Split the single "Database" run into "D" / "atabase" (matches the
# two-run shape produced by Word when text is retyped), using a
# formatting toggle on just the first character to force Word to
# break the run without altering the final look.
$r2 = $d.Content
$r2.Find.Execute("Database")
$firstChar = $r2.Characters(1)
$firstChar.Bold = 1
$firstChar.Bold = 0

# 2. Consolidate "reconditions" + ": " into a single run and drop the
# spell-check proof-error markers around it (Word merges runs with
# identical formatting and clears proofErr tags fully inside the
# replaced span when you re-save over the same text).
$r3 = $d.Content
$r3.Find.Execute("reconditions: ", $false, $false, $false, $false, $false, $true, 1, $false, "reconditions: ", 2)

# 3. Consolidate " (UC-9: " + "Order" + " " + "History" + ")" into a
# single run and drop the proofErr markers around "Order"/"History".
$r4 = $d.Content
$r4.Find.Execute(" (UC-9: Order History)", $false, $false, $false, $false, $false, $true, 1, $false, " (UC-9: Order History)", 2)
